$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 26.81310566666667
$ws.Range("H2").Value = 80.439317
$ws.Range("I2").Value = 0.004518206005002021
$ws.Range("J2").Value = 0.004518206005002021
$ws.Range("M2").Value = 24.576554
$ws.Range("N2").Value = 73.729662
$ws.Range("O2").Value = 0.07553767049546639
$ws.Range("P2").Value = 0.07553767049546638
$ws.Range("Q2").Value = 658.9737393245395
$ws.Range("R2").Value = 5930.763653920854
$ws.Range("S2").Value = 0.0003412947564364802
$ws.Range("T2").Value = 0.0003412947564364801
# Row 3
$ws.Range("G3").Value = 26.81310566666667
$ws.Range("H3").Value = 80.439317
$ws.Range("I3").Value = 0.004518206005002021
$ws.Range("J3").Value = 0.004518206005002021
$ws.Range("O3").Value = 0.359764849016532
$ws.Range("P3").Value = 0.359764849016532
$ws.Range("Q3").Value = 3138.508061989829
$ws.Range("R3").Value = 28246.57255790846
$ws.Range("S3").Value = 0.00162549170121514
$ws.Range("T3").Value = 0.00162549170121514
# Row 4
$ws.Range("G4").Value = 26.81310566666667
$ws.Range("H4").Value = 80.439317
$ws.Range("I4").Value = 0.004518206005002021
$ws.Range("J4").Value = 0.004518206005002021
$ws.Range("M4").Value = 55.68784966666667
$ws.Range("N4").Value = 167.063549
$ws.Range("O4").Value = 0.1711603033819035
$ws.Range("P4").Value = 0.1711603033819035
$ws.Range("Q4").Value = 1493.164197461782
$ws.Range("R4").Value = 13438.47777715603
$ws.Range("S4").Value = 0.0007733375105580842
$ws.Range("T4").Value = 0.0007733375105580842
# Row 5
$ws.Range("G5").Value = 26.81310566666667
$ws.Range("H5").Value = 80.439317
$ws.Range("I5").Value = 0.004518206005002021
$ws.Range("J5").Value = 0.004518206005002021
$ws.Range("M5").Value = 128.0392633333333
$ws.Range("N5").Value = 384.11779
$ws.Range("O5").Value = 0.3935371771060981
$ws.Range("P5").Value = 0.3935371771060981
$ws.Range("Q5").Value = 3433.130297238826
$ws.Range("R5").Value = 30898.17267514943
$ws.Range("S5").Value = 0.001778082036792316
$ws.Range("T5").Value = 0.001778082036792316
# Row 6
$ws.Range("G6").Value = 5771.873535333333
$ws.Range("I6").Value = 0.9726032482643521
$ws.Range("J6").Value = 0.9726032482643523
$ws.Range("M6").Value = 24.576554
$ws.Range("N6").Value = 73.729662
$ws.Range("O6").Value = 0.07553767049546639
$ws.Range("P6").Value = 0.07553767049546638
$ws.Range("Q6").Value = 141852.7616222906
$ws.Range("R6").Value = 1276674.854600615
$ws.Range("S6").Value = 0.07346818369021292
$ws.Range("T6").Value = 0.07346818369021292
# Row 7
$ws.Range("G7").Value = 5771.873535333333
$ws.Range("I7").Value = 0.9726032482643521
$ws.Range("J7").Value = 0.9726032482643523
$ws.Range("O7").Value = 0.359764849016532
$ws.Range("P7").Value = 0.359764849016532
$ws.Range("Q7").Value = 675605.1256661987
$ws.Range("R7").Value = 6080446.130995789
$ws.Range("S7").Value = 0.3499084607648132
$ws.Range("T7").Value = 0.3499084607648133
# Row 8
$ws.Range("G8").Value = 5771.873535333333
$ws.Range("I8").Value = 0.9726032482643521
$ws.Range("J8").Value = 0.9726032482643523
$ws.Range("M8").Value = 55.68784966666667
$ws.Range("N8").Value = 167.063549
$ws.Range("O8").Value = 0.1711603033819035
$ws.Range("P8").Value = 0.1711603033819035
$ws.Range("Q8").Value = 321423.2257306546
$ws.Range("R8").Value = 2892809.031575891
$ws.Range("S8").Value = 0.1664710670431513
$ws.Range("T8").Value = 0.1664710670431514
# Row 9
$ws.Range("G9").Value = 5771.873535333333
$ws.Range("I9").Value = 0.9726032482643521
$ws.Range("J9").Value = 0.9726032482643523
$ws.Range("M9").Value = 128.0392633333333
$ws.Range("N9").Value = 384.11779
$ws.Range("O9").Value = 0.3935371771060981
$ws.Range("P9").Value = 0.3935371771060981
$ws.Range("Q9").Value = 739026.4355172423
$ws.Range("R9").Value = 6651237.919655181
$ws.Range("S9").Value = 0.3827555367661746
$ws.Range("T9").Value = 0.3827555367661747
# Row 10
$ws.Range("G10").Value = 132.4457753333333
$ws.Range("H10").Value = 397.337326
$ws.Range("I10").Value = 0.02231808970163987
$ws.Range("J10").Value = 0.02231808970163988
$ws.Range("M10").Value = 24.576554
$ws.Range("N10").Value = 73.729662
$ws.Range("O10").Value = 0.07553767049546639
$ws.Range("P10").Value = 0.07553767049546638
$ws.Range("Q10").Value = 3255.060749551535
$ws.Range("R10").Value = 29295.54674596382
$ws.Range("S10").Value = 0.001685856505970735
$ws.Range("T10").Value = 0.001685856505970735
# Row 11
$ws.Range("G11").Value = 132.4457753333333
$ws.Range("H11").Value = 397.337326
$ws.Range("I11").Value = 0.02231808970163987
$ws.Range("J11").Value = 0.02231808970163988
$ws.Range("O11").Value = 0.359764849016532
$ws.Range("P11").Value = 0.359764849016532
$ws.Range("Q11").Value = 15502.94616475275
$ws.Range("R11").Value = 139526.5154827748
$ws.Range("S11").Value = 0.008029264171847886
$ws.Range("T11").Value = 0.008029264171847888
# Row 12
$ws.Range("G12").Value = 132.4457753333333
$ws.Range("H12").Value = 397.337326
$ws.Range("I12").Value = 0.02231808970163987
$ws.Range("J12").Value = 0.02231808970163988
$ws.Range("M12").Value = 55.68784966666667
$ws.Range("N12").Value = 167.063549
$ws.Range("O12").Value = 0.1711603033819035
$ws.Range("P12").Value = 0.1711603033819035
$ws.Range("Q12").Value = 7375.620425747776
$ws.Range("R12").Value = 66380.58383172999
$ws.Range("S12").Value = 0.003819971004237217
$ws.Range("T12").Value = 0.003819971004237218
# Row 13
$ws.Range("G13").Value = 132.4457753333333
$ws.Range("H13").Value = 397.337326
$ws.Range("I13").Value = 0.02231808970163987
$ws.Range("J13").Value = 0.02231808970163988
$ws.Range("M13").Value = 128.0392633333333
$ws.Range("N13").Value = 384.11779
$ws.Range("O13").Value = 0.3935371771060981
$ws.Range("P13").Value = 0.3935371771060981
$ws.Range("Q13").Value = 16958.25950529217
$ws.Range("R13").Value = 152624.3355476295
$ws.Range("S13").Value = 0.008782998019584034
$ws.Range("T13").Value = 0.008782998019584036
# Row 14
$ws.Range("G14").Value = 3.326003
$ws.Range("H14").Value = 9.978009
$ws.Range("I14").Value = 0.0005604560290058679
$ws.Range("J14").Value = 0.000560456029005868
$ws.Range("M14").Value = 24.576554
$ws.Range("N14").Value = 73.729662
$ws.Range("O14").Value = 0.07553767049546639
$ws.Range("P14").Value = 0.07553767049546638
$ws.Range("Q14").Value = 81.741692333662
$ws.Range("R14").Value = 735.6752310029581
$ws.Range("S14").Value = 0.0000423355428462428
$ws.Range("T14").Value = 0.0000423355428462428
# Row 15
$ws.Range("G15").Value = 3.326003
$ws.Range("H15").Value = 9.978009
$ws.Range("I15").Value = 0.0005604560290058679
$ws.Range("J15").Value = 0.000560456029005868
$ws.Range("O15").Value = 0.359764849016532
$ws.Range("P15").Value = 0.359764849016532
$ws.Range("Q15").Value = 389.3128740651424
$ws.Range("R15").Value = 3503.815866586281
$ws.Range("S15").Value = 0.0002016323786557011
$ws.Range("T15").Value = 0.0002016323786557012
# Row 16
$ws.Range("G16").Value = 3.326003
$ws.Range("H16").Value = 9.978009
$ws.Range("I16").Value = 0.0005604560290058679
$ws.Range("J16").Value = 0.000560456029005868
$ws.Range("M16").Value = 55.68784966666667
$ws.Range("N16").Value = 167.063549
$ws.Range("O16").Value = 0.1711603033819035
$ws.Range("P16").Value = 0.1711603033819035
$ws.Range("Q16").Value = 185.2179550548824
$ws.Range("R16").Value = 1666.961595493941
$ws.Range("S16").Value = 0.00009592782395686127
$ws.Range("T16").Value = 0.00009592782395686128
# Row 17
$ws.Range("G17").Value = 3.326003
$ws.Range("H17").Value = 9.978009
$ws.Range("I17").Value = 0.0005604560290058679
$ws.Range("J17").Value = 0.000560456029005868
$ws.Range("M17").Value = 128.0392633333333
$ws.Range("N17").Value = 384.11779
$ws.Range("O17").Value = 0.3935371771060981
$ws.Range("P17").Value = 0.3935371771060981
$ws.Range("Q17").Value = 425.8589739644567
$ws.Range("R17").Value = 3832.73076568011
$ws.Range("S17").Value = 0.0002205602835470627
$ws.Range("T17").Value = 0.0002205602835470627
